$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.237.57"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.661.65"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06337"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "1.662.74"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "1.890.39"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5473"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "0.0₅8215"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "26.272.61"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.652"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.087"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1240"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.243"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.417"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05945"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.646"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.306"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.631"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.785"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5896"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01592"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.949"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8599"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "1.027.77"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "1.803.20"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  +7.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.096"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05187"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.468"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.15%  "
